$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = "10/05/2025"
$ws.Range("A8").Style = "Normal"

$ws.Range("B8").Value = 2602100114

$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "14:06:51"
$ws.Range("C8").Style = "Normal"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = ""
$ws.Range("D8").Style = "Normal"

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = ""
$ws.Range("E8").Style = "Normal"
